$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New string order must match: B14, C14, E14, B15, C15, D15, E15, F14, G14, F15, G15
$ws.Range("B14").Value = "Doc Dash See appointment list"
$ws.Range("C14").Value = "the doctor can see the appointment list only to the given doctor"
$ws.Range("E14").Value = "http://localhost:8080/appointments/docapp/{doctorID}"

$ws.Range("B15").Value = "Update Api to Approve the appointment"
$ws.Range("C15").Value = "api used to approve the appointment"
$ws.Range("D15").Value = "PUT"
$ws.Range("E15").Value = "http://localhost:8080/appointments/approve/{appointment ID}"

$ws.Range("F14").Value = "DoctorID"
$ws.Range("G14").Value = "all appointment details"

$ws.Range("F15").Value = "Appointment ID"
$ws.Range("G15").Value = "change the state from pending to approved "

# D14 reuses the existing "GET" shared string
$ws.Range("D14").Value = "GET"

# Hyperlinks for the URL cells
$ws.Hyperlinks.Add($ws.Range("E14"), "http://localhost:8080/appointments/docapp/{doctorID}")
$ws.Range("E14").HorizontalAlignment = -4108
$ws.Hyperlinks.Add($ws.Range("E15"), "http://localhost:8080/appointments/approve/{appointment ID}")
$ws.Range("E15").HorizontalAlignment = -4108

# Update view: scroll/selection
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("G19").Select()
